$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Volue load-forecast data (3 days later); values shifted per upstream refetch.
$timestamps = @(
    45427,
    45427.01041666666,
    45427.02083333334,
    45427.03125,
    45427.04166666666,
    45427.05208333334,
    45427.0625,
    45427.07291666666,
    45427.08333333334,
    45427.09375,
    45427.10416666666,
    45427.11458333334,
    45427.125,
    45427.13541666666,
    45427.14583333334,
    45427.15625,
    45427.16666666666,
    45427.17708333334,
    45427.1875,
    45427.19791666666,
    45427.20833333334,
    45427.21875,
    45427.22916666666,
    45427.23958333334,
    45427.25,
    45427.26041666666,
    45427.27083333334,
    45427.28125,
    45427.29166666666,
    45427.30208333334,
    45427.3125,
    45427.32291666666,
    45427.33333333334,
    45427.34375,
    45427.35416666666,
    45427.36458333334,
    45427.375,
    45427.38541666666,
    45427.39583333334,
    45427.40625,
    45427.41666666666,
    45427.42708333334,
    45427.4375,
    45427.44791666666,
    45427.45833333334,
    45427.46875,
    45427.47916666666,
    45427.48958333334,
    45427.5,
    45427.51041666666,
    45427.52083333334,
    45427.53125,
    45427.54166666666,
    45427.55208333334,
    45427.5625,
    45427.57291666666,
    45427.58333333334,
    45427.59375,
    45427.60416666666,
    45427.61458333334,
    45427.625,
    45427.63541666666,
    45427.64583333334,
    45427.65625,
    45427.66666666666,
    45427.67708333334,
    45427.6875,
    45427.69791666666,
    45427.70833333334,
    45427.71875,
    45427.72916666666,
    45427.73958333334,
    45427.75,
    45427.76041666666,
    45427.77083333334,
    45427.78125,
    45427.79166666666,
    45427.80208333334,
    45427.8125,
    45427.82291666666,
    45427.83333333334,
    45427.84375,
    45427.85416666666,
    45427.86458333334,
    45427.875,
    45427.88541666666,
    45427.89583333334,
    45427.90625,
    45427.91666666666,
    45427.92708333334,
    45427.9375,
    45427.94791666666,
    45427.95833333334,
    45427.96875,
    45427.97916666666,
    45427.98958333334,
    45428
)
$loads = @(
    5380,
    5340,
    5300,
    5260,
    5230,
    5210,
    5190,
    5160,
    5120,
    5100,
    5100,
    5110,
    5130,
    5150,
    5160,
    5170,
    5190,
    5220,
    5260,
    5330,
    5420,
    5520,
    5620,
    5750,
    5900,
    6030,
    6110,
    6160,
    6180,
    6180,
    6170,
    6140,
    6060,
    5970,
    5880,
    5780,
    5710,
    5630,
    5560,
    5500,
    5440,
    5420,
    5400,
    5390,
    5380,
    5380,
    5380,
    5380,
    5390,
    5400,
    5400,
    5380,
    5340,
    5320,
    5310,
    5310,
    5310,
    5310,
    5320,
    5340,
    5390,
    5450,
    5500,
    5560,
    5610,
    5660,
    5700,
    5770,
    5890,
    6000,
    6110,
    6220,
    6320,
    6410,
    6500,
    6600,
    6700,
    6810,
    6910,
    7040,
    7130,
    7130,
    7120,
    7080,
    6930,
    6720,
    6480,
    6300,
    6140,
    6020,
    5910,
    5800,
    5720,
    5640,
    5570,
    5520,
    5470
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $timestamps[$i]
    $ws.Cells.Item($row, 2).Value = $loads[$i]
}
